$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
"before ScrollRow: " + $win.ScrollRow
$win.ScrollRow = 9
$win.ScrollColumn = 1
"after ScrollRow: " + $win.ScrollRow
"after ScrollColumn: " + $win.ScrollColumn
$ws.Range("I14").Select()
"done"
